$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '68.843.06'
$ws.Range('E2').Value = '  +1.36%  '
$ws.Range('D3').Value = '3.867.77'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '602.34'
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('D6').Value = '172.73'
$ws.Range('E6').Value = '  +3.61%  '
$ws.Range('D7').Value = '3.868.53'
$ws.Range('E7').Value = '  +0.74%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +0.70%  '
$ws.Range('D10').Value = '0.169'
$ws.Range('E10').Value = '  +2.43%  '
$ws.Range('D11').Value = '6.55'
$ws.Range('E11').Value = '  +3.86%  '
$ws.Range('D12').Value = '0.463'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('D13').Value = '0.0000287'
$ws.Range('E13').Value = '  +15.73%  '
$ws.Range('D14').Value = '37.37'
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('D15').Value = '4.505.59'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '3.886.66'
$ws.Range('E16').Value = '  +0.83%  '
$ws.Range('D17').Value = '68.844.09'
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').Value = '18.44'
$ws.Range('E18').Value = '  +1.75%  '
$ws.Range('D19').Value = '7.41'
$ws.Range('E19').Value = '  -0.86%  '
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('D21').Value = '11.21'
$ws.Range('E21').Value = '  +4.14%  '
$ws.Range('D22').Value = '474.75'
$ws.Range('E22').Value = '  +1.04%  '
$ws.Range('D23').Value = '0.733'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').Value = '0.0000163'
$ws.Range('E24').Value = '  +1.31%  '
$ws.Range('D25').Value = '83.93'
$ws.Range('E25').Value = '  -0.21%  '
$ws.Range('D26').Value = '2.28'
$ws.Range('E26').Value = '  +2.87%  '
$ws.Range('D27').Value = '12.20'
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').Value = '10.51'
$ws.Range('E28').Value = '  +5.15%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').Value = '2.95'
$ws.Range('E30').Value = '  +1.18%  '
$ws.Range('D31').Value = '4.007.80'
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('D32').Value = '7.81'
$ws.Range('E32').Value = '  +1.20%  '
$ws.Range('D33').Value = '31.57'
$ws.Range('E33').Value = '  +1.81%  '
$ws.Range('D34').Value = '2.32'
$ws.Range('E34').Value = '  +0.64%  '
$ws.Range('D35').Value = '9.43'
$ws.Range('E35').Value = '  +0.68%  '
$ws.Range('D36').Value = '3.824.44'
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('D37').Value = '4.03'
$ws.Range('E37').Value = '  +22.93%  '
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('E39').Value = '  +0.82%  '
$ws.Range('E40').Value = '  +1.17%  '
$ws.Range('D41').Value = '5.98'
$ws.Range('E41').Value = '  +1.05%  '
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').Value = '0.322'
$ws.Range('E43').Value = '  +2.77%  '
$ws.Range('D44').Value = '2.01'
$ws.Range('E44').Value = '  +0.76%  '
$ws.Range('D45').Value = '0.000301'
$ws.Range('E45').Value = '  +11.23%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = '424.78'
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('B47').Value = 'USDe'
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('D48').Value = '8.72'
$ws.Range('E48').Value = '  +1.66%  '
$ws.Range('D49').Value = '46.46'
$ws.Range('E49').Value = '  -1.99%  '
$ws.Range('D50').Value = '142.14'
$ws.Range('E50').Value = '  -0.61%  '
$ws.Range('D51').Value = '0.0360'
$ws.Range('E51').Value = '  +0.94%  '
